# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new F value for each of the affected sheets
$updates = @{
    2  = 313
    4  = 10334
    7  = 1282
    8  = 7127
    13 = 3183
    18 = 958
    20 = 64
    21 = 1626
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
